$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "E2" = 3
    "F2" = 1
    "G2" = 76.16539233333333
    "H2" = 228.496177
    "I2" = 0.5742142031125765
    "J2" = 0.5742142031125764
    "K2" = 3
    "L2" = 1
    "M2" = 3.771625333333333
    "N2" = 11.314876
    "O2" = 0.7452169260344172
    "P2" = 0.7452169260344171
    "Q2" = 287.2673232476724
    "R2" = 2585.405909229052
    "S2" = 0.4279141433288568
    "T2" = 0.4279141433288566
    "E3" = 3
    "F3" = 1
    "G3" = 76.16539233333333
    "H3" = 228.496177
    "I3" = 0.5742142031125765
    "J3" = 0.5742142031125764
    "M3" = 0.3326733333333333
    "N3" = 0.99802
    "O3" = 0.06573129007519563
    "P3" = 0.06573129007519561
    "Q3" = 25.33819495217111
    "R3" = 228.04375456954
    "S3" = 0.03774384035009007
    "T3" = 0.03774384035009005
    "E4" = 3
    "F4" = 1
    "G4" = 76.16539233333333
    "H4" = 228.496177
    "I4" = 0.5742142031125765
    "J4" = 0.5742142031125764
    "M4" = 0.3774316666666667
    "N4" = 1.132295
    "O4" = 0.07457486933698085
    "P4" = 0.07457486933698083
    "Q4" = 28.74723097069056
    "R4" = 258.725078736215
    "S4" = 0.04282194916855898
    "T4" = 0.04282194916855896
    "E5" = 3
    "F5" = 1
    "G5" = 76.16539233333333
    "H5" = 228.496177
    "I5" = 0.5742142031125765
    "J5" = 0.5742142031125764
    "K5" = 3
    "L5" = 1
    "M5" = 0.5793803333333334
    "N5" = 1.738141
    "O5" = 0.1144769145534063
    "P5" = 0.1144769145534063
    "Q5" = 44.12873039855078
    "R5" = 397.158573586957
    "S5" = 0.06573427026507074
    "T5" = 0.06573427026507071
    "I6" = 0.03912478832313545
    "J6" = 0.03912478832313544
    "K6" = 3
    "L6" = 1
    "M6" = 3.771625333333333
    "N6" = 11.314876
    "O6" = 0.7452169260344172
    "P6" = 0.7452169260344171
    "Q6" = 19.57331106283244
    "R6" = 176.159799565492
    "S6" = 0.02915645448591426
    "T6" = 0.02915645448591425
    "I7" = 0.03912478832313545
    "J7" = 0.03912478832313544
    "M7" = 0.3326733333333333
    "N7" = 0.99802
    "O7" = 0.06573129007519563
    "P7" = 0.06573129007519561
    "S7" = 0.002571722810398643
    "T7" = 0.002571722810398642
    "I8" = 0.03912478832313545
    "J8" = 0.03912478832313544
    "M8" = 0.3774316666666667
    "N8" = 1.132295
    "O8" = 0.07457486933698085
    "P8" = 0.07457486933698083
    "Q8" = 1.958727806640556
    "R8" = 17.628550259765
    "S8" = 0.00291772597703486
    "T8" = 0.002917725977034859
    "I9" = 0.03912478832313545
    "J9" = 0.03912478832313544
    "K9" = 3
    "L9" = 1
    "M9" = 0.5793803333333334
    "N9" = 1.738141
    "O9" = 0.1144769145534063
    "P9" = 0.1144769145534063
    "Q9" = 3.006765117360778
    "R9" = 27.060886056247
    "S9" = 0.004478885049787687
    "T9" = 0.004478885049787685
    "G10" = 51.18420533333333
    "H10" = 153.552616
    "I10" = 0.3858799485835225
    "J10" = 0.3858799485835225
    "K10" = 3
    "L10" = 1
    "M10" = 3.771625333333333
    "N10" = 11.314876
    "O10" = 0.7452169260344172
    "P10" = 0.7452169260344171
    "Q10" = 193.0476455017351
    "R10" = 1737.428809515616
    "S10" = 0.2875642691017317
    "T10" = 0.2875642691017316
    "G11" = 51.18420533333333
    "H11" = 153.552616
    "I11" = 0.3858799485835225
    "J11" = 0.3858799485835225
    "M11" = 0.3326733333333333
    "N11" = 0.99802
    "O11" = 0.06573129007519563
    "P11" = 0.06573129007519561
    "Q11" = 17.02762020225778
    "R11" = 153.24858182032
    "S11" = 0.02536438683454509
    "T11" = 0.02536438683454509
    "G12" = 51.18420533333333
    "H12" = 153.552616
    "I12" = 0.3858799485835225
    "J12" = 0.3858799485835225
    "M12" = 0.3774316666666667
    "N12" = 1.132295
    "O12" = 0.07457486933698085
    "P12" = 0.07457486933698083
    "Q12" = 19.31853992596889
    "R12" = 173.86685933372
    "S12" = 0.02877694674537708
    "T12" = 0.02877694674537708
    "G13" = 51.18420533333333
    "H13" = 153.552616
    "I13" = 0.3858799485835225
    "J13" = 0.3858799485835225
    "K13" = 3
    "L13" = 1
    "M13" = 0.5793803333333334
    "N13" = 1.738141
    "O13" = 0.1144769145534063
    "P13" = 0.1144769145534063
    "Q13" = 29.65512194742844
    "R13" = 266.896097526856
    "S13" = 0.04417434590186874
    "T13" = 0.04417434590186874
    "E14" = 2
    "F14" = 0.6666666666666666
    "G14" = 0.103602
    "H14" = 0.310806
    "I14" = 0.0007810599807654878
    "J14" = 0.0007810599807654877
    "K14" = 3
    "L14" = 1
    "M14" = 3.771625333333333
    "N14" = 11.314876
    "O14" = 0.7452169260344172
    "P14" = 0.7452169260344171
    "Q14" = 0.390747927784
    "R14" = 3.516731350056
    "S14" = 0.0005820591179145579
    "T14" = 0.0005820591179145577
    "E15" = 2
    "F15" = 0.6666666666666666
    "G15" = 0.103602
    "H15" = 0.310806
    "I15" = 0.0007810599807654878
    "J15" = 0.0007810599807654877
    "M15" = 0.3326733333333333
    "N15" = 0.99802
    "O15" = 0.06573129007519563
    "P15" = 0.06573129007519561
    "Q15" = 0.03446562268
    "R15" = 0.31019060412
    "S15" = [double]"5.134008016182299E-05"
    "T15" = [double]"5.134008016182297E-05"
    "E16" = 2
    "F16" = 0.6666666666666666
    "G16" = 0.103602
    "H16" = 0.310806
    "I16" = 0.0007810599807654878
    "J16" = 0.0007810599807654877
    "M16" = 0.3774316666666667
    "N16" = 1.132295
    "O16" = 0.07457486933698085
    "P16" = 0.07457486933698083
    "Q16" = 0.03910267553000001
    "R16" = 0.3519240797700001
    "S16" = [double]"5.824744600993103E-05"
    "T16" = [double]"5.824744600993101E-05"
    "E17" = 2
    "F17" = 0.6666666666666666
    "G17" = 0.103602
    "H17" = 0.310806
    "I17" = 0.0007810599807654878
    "J17" = 0.0007810599807654877
    "K17" = 3
    "L17" = 1
    "M17" = 0.5793803333333334
    "N17" = 1.738141
    "O17" = 0.1144769145534063
    "P17" = 0.1144769145534063
    "Q17" = 0.06002496129400001
    "R17" = 0.5402246516460001
    "S17" = [double]"8.941333667917595E-05"
    "T17" = [double]"8.941333667917592E-05"
}

foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}

Write-Output "Applied $($changes.Count) cell updates"
